$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: create row 29 as a new data row, copying the formatting of row 28
# (only columns A:R are used by the sheet, so restrict the copy/paste range
# to avoid touching the full 16384-column row and bloating the dimension).
$ws.Range("A28:R28").Copy()
$ws.Range("A29:R29").PasteSpecial(-4122)  # xlPasteFormats

# Populate the constant (non-varying) columns of the new row 29 with the
# same values used throughout this block of rows.
$ws.Range("A29").Value = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 100112043
$ws.Range("G29").Value = "Pepino dulce"
$ws.Range("H29").Value = "Cultivar IV Región"
$ws.Range("N29").Value = "`$/bandeja 18 kilos"
$ws.Range("O29").Value = "Provincia de Limarí"
$ws.Range("Q29").Value = 18
$ws.Range("R29").Value = "Hortaliza"

# --- Step 2: update the varying columns (D, I, J, K, L, M, P) for rows 20-29.
# Row 20
$ws.Range("D20").Value = 44719
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 15000
$ws.Range("P20").Value = 833
# Row 21
$ws.Range("D21").Value = 44396
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 15000
$ws.Range("P21").Value = 833
# Row 22
$ws.Range("D22").Value = 44396
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 12000
$ws.Range("L22").Value = 12000
$ws.Range("M22").Value = 12000
$ws.Range("P22").Value = 667
# Row 23
$ws.Range("D23").Value = 44596
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 150
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 14000
$ws.Range("M23").Value = 14000
$ws.Range("P23").Value = 778
# Row 24
$ws.Range("D24").Value = 44630
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 15000
$ws.Range("P24").Value = 833
# Row 25
$ws.Range("D25").Value = 44627
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 15000
$ws.Range("P25").Value = 833
# Row 26
$ws.Range("D26").Value = 44245
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 12000
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = 12000
$ws.Range("P26").Value = 667
# Row 27
$ws.Range("D27").Value = 44245
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = 10000
$ws.Range("P27").Value = 556
# Row 28
$ws.Range("D28").Value = 44249
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 400
$ws.Range("K28").Value = 12000
$ws.Range("L28").Value = 12000
$ws.Range("M28").Value = 12000
$ws.Range("P28").Value = 667
# Row 29
$ws.Range("D29").Value = 44249
$ws.Range("I29").Value = "Segunda"
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = 10000
$ws.Range("P29").Value = 556
